$p = $ppt.ActivePresentation

# The deck's (slide-master) theme was re-coloured from the custom "Integral"
# palette to the stock PowerPoint "Office Theme" palette (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink - the 12 slots of <a:clrScheme>), matching the
# colours PowerPoint ships as its default "Office Theme".
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$tcs = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
